# Fixed update to excel issue
#
# 1. Rename "Requested quantity" header on "Weekly Quantity" -> "Weekly_PO_Qty"
# 2. Rename "Requested quantity" header on "Monthly Trend"   -> "Monthly_PO_Qty"
# 3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with
#    ds / PO_Forecast / yhat_lower / yhat_upper columns and 49 rows of data.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: fix up existing headers ------------------------------------
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the new "PO Forecast" sheet after "Monthly Trend" ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the header styling (bold font, thin border, center/top alignment)
# and the date-column styling already used on the other two sheets, so we
# don't introduce brand-new style entries.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A50").PasteSpecial(-4122)

# Data rows (r=2..50)
$poForecastData = @(
    @(45067.99999999999, 99, -4.396048961306656, 210.233596269684),
    @(45081.99999999999, 99, -20.2165661089829, 210.1436766058947),
    @(45088.99999999999, 99, -12.91427062644268, 212.1566434817589),
    @(45095.99999999999, 98, -6.489327121757212, 206.5586176303883),
    @(45102.99999999999, 98, -11.17751658460789, 210.930910087773),
    @(45109.99999999999, 98, -17.95987483454107, 206.0882431062323),
    @(45116.99999999999, 97, -2.950134337579711, 203.0386008490146),
    @(45123.99999999999, 97, -12.69651905068577, 208.7540369775691),
    @(45130.99999999999, 97, -11.73661982148063, 205.5446089396035),
    @(45144.99999999999, 96, -12.74761416630265, 199.8682168527994),
    @(45151.99999999999, 96, -6.37761584665612, 217.7593522640513),
    @(45158.99999999999, 96, -17.06744268315429, 207.4707369057402),
    @(45172.99999999999, 95, -13.44007698497807, 197.7972252995544),
    @(45179.99999999999, 95, -17.00611810487149, 199.9128728494959),
    @(45186.99999999999, 95, -19.49845835351579, 206.9146262279759),
    @(45200.99999999999, 94, -8.287952603374402, 206.6193788305959),
    @(45207.99999999999, 94, -11.31977234325064, 209.5488005254486),
    @(45214.99999999999, 93, -17.63676687346891, 202.6697336831358),
    @(45221.99999999999, 93, -20.44747371522464, 196.5400943399612),
    @(45228.99999999999, 93, -10.59487456679438, 199.7157251945999),
    @(45235.99999999999, 93, -14.1120126857293, 199.1856384831627),
    @(45242.99999999999, 92, -17.07332754966603, 196.9092421163638),
    @(45249.99999999999, 92, -15.06931977807919, 201.5131341880353),
    @(45256.99999999999, 92, -14.78781884186865, 192.411884374235),
    @(45270.99999999999, 91, -10.11559907110237, 207.8718075589789),
    @(45277.99999999999, 91, -15.46613280070124, 202.1967001874025),
    @(45298.99999999999, 90, -18.74919951592252, 190.1100854315006),
    @(45305.99999999999, 90, -24.4349614067728, 197.2536190491094),
    @(45312.99999999999, 89, -24.41585948624586, 201.2633221841409),
    @(45319.99999999999, 89, -10.88077617133848, 198.9985983612795),
    @(45326.99999999999, 89, -16.46639077853176, 197.1902115893416),
    @(45333.99999999999, 89, -22.64900996096026, 197.3108525415806),
    @(45347.99999999999, 88, -22.17224265399865, 193.3779895030701),
    @(45354.99999999999, 88, -21.75638115979671, 196.818029210622),
    @(45361.99999999999, 87, -24.19301122817012, 189.1216434656662),
    @(45375.99999999999, 87, -23.17530506724622, 194.9213675809269),
    @(45389.99999999999, 86, -27.89985362271038, 193.773050961222),
    @(45403.99999999999, 86, -29.10395646118725, 198.1311652135032),
    @(45410.99999999999, 85, -26.3626315073152, 192.8337430469729),
    @(45417.99999999999, 85, -26.60510587580034, 193.6566813266842),
    @(45424.99999999999, 85, -18.77397123678445, 192.352129314073),
    @(45431.99999999999, 84, -20.83721328175757, 204.3203830217684),
    @(45438.99999999999, 84, -23.52499340176914, 192.8591376226714),
    @(45445.99999999999, 84, -24.50981537511574, 201.5782918269426),
    @(45452.99999999999, 84, -24.53975878470312, 191.8066767888909),
    @(45459.99999999999, 83, -22.81074051617579, 192.0749670597956),
    @(45466.99999999999, 83, -26.0391182146705, 198.0285053791373),
    @(45473.99999999999, 83, -31.15590175155913, 189.3513472624174),
    @(45480.99999999999, 82, -26.34415934342852, 191.4640816498016)
)

$r = 2
foreach ($row in $poForecastData) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}
